$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (target OOXML widths are 15.42578125 / 14.7109375,
# i.e. exactly 1 character narrower than the original 16.42578125 / 15.7109375).
# The COM ColumnWidth setter here quantizes to a 1/6-character grid, so we
# feed it the value whose rounded result lands closest to the target width.
$ws.Columns.Item(1).ColumnWidth = 14.592447916666666
$ws.Columns.Item(2).ColumnWidth = 13.877604166666666

# Update cell values in A1:B32
$ws.Cells.Item(1, 1).Value = -0.2267693919272773
$ws.Cells.Item(1, 2).Value = 0.2265989915284905
$ws.Cells.Item(2, 1).Value = -0.20449639056043889
$ws.Cells.Item(2, 2).Value = 0.20380148102344897
$ws.Cells.Item(3, 1).Value = -0.10382851170790452
$ws.Cells.Item(3, 2).Value = 0.10355405099605619
$ws.Cells.Item(4, 1).Value = -0.095554051092774372
$ws.Cells.Item(4, 2).Value = 0.09508730421719136
$ws.Cells.Item(5, 1).Value = -0.092087304273766435
$ws.Cells.Item(5, 2).Value = 0.090497240551423275
$ws.Cells.Item(6, 1).Value = -0.05267360220647177
$ws.Cells.Item(6, 2).Value = 0.052078648558984142
$ws.Cells.Item(7, 1).Value = -0.042078648696300291
$ws.Cells.Item(7, 2).Value = 0.041928825570896944
$ws.Cells.Item(8, 1).Value = -0.031928825712881803
$ws.Cells.Item(8, 2).Value = 0.031646907569003346
$ws.Cells.Item(9, 1).Value = -0.029646907643547937
$ws.Cells.Item(9, 2).Value = 0.029410087352455427
$ws.Cells.Item(10, 1).Value = -0.027410087432551578
$ws.Cells.Item(10, 2).Value = 0.027395181733666618
$ws.Cells.Item(11, 1).Value = -0.024395181823551937
$ws.Cells.Item(11, 2).Value = 0.024369043809220337
$ws.Cells.Item(12, 1).Value = -0.02086904390496791
$ws.Cells.Item(12, 2).Value = 0.020672028952176014
$ws.Cells.Item(13, 1).Value = -0.017172029053683424
$ws.Cells.Item(13, 2).Value = 0.017082623813123199
$ws.Cells.Item(14, 1).Value = -0.0090826239575259038
$ws.Cells.Item(14, 2).Value = 0.0090537210597760875
$ws.Cells.Item(15, 1).Value = -0.0080537211411879639
$ws.Cells.Item(15, 2).Value = 0.0080350807127045343
$ws.Cells.Item(16, 1).Value = -0.0060350808042910487
$ws.Cells.Item(16, 2).Value = 0.0060033072695362399
$ws.Cells.Item(17, 1).Value = -0.0040033073626464244
$ws.Cells.Item(17, 2).Value = 0.0039999998881734555
$ws.Cells.Item(18, 1).Value = -0.016102601025302477
$ws.Cells.Item(18, 2).Value = 0.016090935491671843
$ws.Cells.Item(19, 1).Value = -0.012090935530614466
$ws.Cells.Item(19, 2).Value = 0.012016130962197114
$ws.Cells.Item(20, 1).Value = -0.0080161310038899813
$ws.Cells.Item(20, 2).Value = 0.008005619013660592
$ws.Cells.Item(21, 1).Value = -0.0040056190557846705
$ws.Cells.Item(21, 2).Value = 0.0039999999575321965
$ws.Cells.Item(22, 1).Value = -0.0457055898658556
$ws.Cells.Item(22, 2).Value = 0.045494536240504146
$ws.Cells.Item(23, 1).Value = -0.040494536307243756
$ws.Cells.Item(23, 2).Value = 0.040098028044343259
$ws.Cells.Item(24, 1).Value = -0.020098028253709543
$ws.Cells.Item(24, 2).Value = 0.019999999787942535
$ws.Cells.Item(25, 1).Value = -0.035823638407117997
$ws.Cells.Item(25, 2).Value = 0.035803356896359162
$ws.Cells.Item(26, 1).Value = -0.033303356963793718
$ws.Cells.Item(26, 2).Value = 0.033279532295669156
$ws.Cells.Item(27, 1).Value = -0.030779532364270334
$ws.Cells.Item(27, 2).Value = 0.030648483832369067
$ws.Cells.Item(28, 1).Value = -0.028648483901972277
$ws.Cells.Item(28, 2).Value = 0.028573459868383111
$ws.Cells.Item(29, 1).Value = -0.021573459987528132
$ws.Cells.Item(29, 2).Value = 0.021561937789405405
$ws.Cells.Item(30, 1).Value = 0.038438061605584384
$ws.Cells.Item(30, 2).Value = -0.038533447554851552
$ws.Cells.Item(31, 1).Value = 0.037820935324427651
$ws.Cells.Item(31, 2).Value = -0.037861087317322628
$ws.Cells.Item(32, 1).Value = 0.047861087176654493
$ws.Cells.Item(32, 2).Value = -0.047993854734903252
